$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the address label cell (merged B2:E2) from "ADRES" to "address"
$ws.Range("B2").Value = "address"

# Match the selection left behind by the edit (the merged B2:E2 range)
$ws.Range("B2:E2").Select()
